$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Resolving-Mac" target cluster was dropped from the upstream NATMI
# TPM re-run, so every remaining row's numbers were recomputed against the
# new (3-cluster) input. Trigger a shared-string-table rebuild while the
# soon-to-be-deleted rows still reference "Resolving-Mac" so the saved
# workbook doesn't keep an orphaned <si> entry around.
$ws.Range("A1:T13").Replace("Resolving-Mac", "MuSCs") | Out-Null

# Drop the three rows that sent to the now-absent "Resolving-Mac" cluster.
$ws.Rows("11:13").Delete()


# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Vegfc"
$ws.Cells.Item(2,3).Value = "Flt4"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 5.441829000000001
$ws.Cells.Item(2,8).Value = 16.325487
$ws.Cells.Item(2,9).Value = 0.5729403216841985
$ws.Cells.Item(2,10).Value = 0.5729403216841985
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 16.75848066666667
$ws.Cells.Item(2,14).Value = 50.275442
$ws.Cells.Item(2,15).Value = 0.97993745062104
$ws.Cells.Item(2,16).Value = 0.9799374506210401
$ws.Cells.Item(2,17).Value = 91.19678608780602
$ws.Cells.Item(2,18).Value = 820.7710747902541
$ws.Cells.Item(2,19).Value = 0.561445678189212
$ws.Cells.Item(2,20).Value = 0.5614456781892121

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Vegfc"
$ws.Cells.Item(3,3).Value = "Flt4"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 5.441829000000001
$ws.Cells.Item(3,8).Value = 16.325487
$ws.Cells.Item(3,9).Value = 0.5729403216841985
$ws.Cells.Item(3,10).Value = 0.5729403216841985
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.1099223333333333
$ws.Cells.Item(3,14).Value = 0.329767
$ws.Cells.Item(3,15).Value = 0.006427611979601263
$ws.Cells.Item(3,16).Value = 0.006427611979601264
$ws.Cells.Item(3,17).Value = 0.5981785412810001
$ws.Cells.Item(3,18).Value = 5.383606871529
$ws.Cells.Item(3,19).Value = 0.003682638075253955
$ws.Cells.Item(3,20).Value = 0.003682638075253956

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Vegfc"
$ws.Cells.Item(4,3).Value = "Flt4"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 5.441829000000001
$ws.Cells.Item(4,8).Value = 16.325487
$ws.Cells.Item(4,9).Value = 0.5729403216841985
$ws.Cells.Item(4,10).Value = 0.5729403216841985
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.233179
$ws.Cells.Item(4,14).Value = 0.6995370000000001
$ws.Cells.Item(4,15).Value = 0.01363493739935873
$ws.Cells.Item(4,16).Value = 0.01363493739935873
$ws.Cells.Item(4,17).Value = 1.268920244391
$ws.Cells.Item(4,18).Value = 11.420282199519
$ws.Cells.Item(4,19).Value = 0.007812005419732497
$ws.Cells.Item(4,20).Value = 0.007812005419732499

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Vegfc"
$ws.Cells.Item(5,3).Value = "Flt4"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.579868
$ws.Cells.Item(5,8).Value = 7.739604
$ws.Cells.Item(5,9).Value = 0.2716201486343598
$ws.Cells.Item(5,10).Value = 0.2716201486343598
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 16.75848066666667
$ws.Cells.Item(5,14).Value = 50.275442
$ws.Cells.Item(5,15).Value = 0.97993745062104
$ws.Cells.Item(5,16).Value = 0.9799374506210401
$ws.Cells.Item(5,17).Value = 43.234668000552
$ws.Cells.Item(5,18).Value = 389.112012004968
$ws.Cells.Item(5,19).Value = 0.2661707559900625
$ws.Cells.Item(5,20).Value = 0.2661707559900625

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Vegfc"
$ws.Cells.Item(6,3).Value = "Flt4"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.579868
$ws.Cells.Item(6,8).Value = 7.739604
$ws.Cells.Item(6,9).Value = 0.2716201486343598
$ws.Cells.Item(6,10).Value = 0.2716201486343598
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.1099223333333333
$ws.Cells.Item(6,14).Value = 0.329767
$ws.Cells.Item(6,15).Value = 0.006427611979601263
$ws.Cells.Item(6,16).Value = 0.006427611979601264
$ws.Cells.Item(6,17).Value = 0.283585110252
$ws.Cells.Item(6,18).Value = 2.552265992268
$ws.Cells.Item(6,19).Value = 0.001745868921263287
$ws.Cells.Item(6,20).Value = 0.001745868921263287

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Vegfc"
$ws.Cells.Item(7,3).Value = "Flt4"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.579868
$ws.Cells.Item(7,8).Value = 7.739604
$ws.Cells.Item(7,9).Value = 0.2716201486343598
$ws.Cells.Item(7,10).Value = 0.2716201486343598
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.233179
$ws.Cells.Item(7,14).Value = 0.6995370000000001
$ws.Cells.Item(7,15).Value = 0.01363493739935873
$ws.Cells.Item(7,16).Value = 0.01363493739935873
$ws.Cells.Item(7,17).Value = 0.601571040372
$ws.Cells.Item(7,18).Value = 5.414139363348
$ws.Cells.Item(7,19).Value = 0.003703523723034008
$ws.Cells.Item(7,20).Value = 0.003703523723034009

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Vegfc"
$ws.Cells.Item(8,3).Value = "Flt4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.476376
$ws.Cells.Item(8,8).Value = 4.429128
$ws.Cells.Item(8,9).Value = 0.1554395296814417
$ws.Cells.Item(8,10).Value = 0.1554395296814417
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 16.75848066666667
$ws.Cells.Item(8,14).Value = 50.275442
$ws.Cells.Item(8,15).Value = 0.97993745062104
$ws.Cells.Item(8,16).Value = 0.9799374506210401
$ws.Cells.Item(8,17).Value = 24.74181865273067
$ws.Cells.Item(8,18).Value = 222.676367874576
$ws.Cells.Item(8,19).Value = 0.1523210164417654
$ws.Cells.Item(8,20).Value = 0.1523210164417655

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Vegfc"
$ws.Cells.Item(9,3).Value = "Flt4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.476376
$ws.Cells.Item(9,8).Value = 4.429128
$ws.Cells.Item(9,9).Value = 0.1554395296814417
$ws.Cells.Item(9,10).Value = 0.1554395296814417
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.1099223333333333
$ws.Cells.Item(9,14).Value = 0.329767
$ws.Cells.Item(9,15).Value = 0.006427611979601263
$ws.Cells.Item(9,16).Value = 0.006427611979601264
$ws.Cells.Item(9,17).Value = 0.1622866947973333
$ws.Cells.Item(9,18).Value = 1.460580253176
$ws.Cells.Item(9,19).Value = 0.0009991049830840204
$ws.Cells.Item(9,20).Value = 0.0009991049830840209

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Vegfc"
$ws.Cells.Item(10,3).Value = "Flt4"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.476376
$ws.Cells.Item(10,8).Value = 4.429128
$ws.Cells.Item(10,9).Value = 0.1554395296814417
$ws.Cells.Item(10,10).Value = 0.1554395296814417
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.233179
$ws.Cells.Item(10,14).Value = 0.6995370000000001
$ws.Cells.Item(10,15).Value = 0.01363493739935873
$ws.Cells.Item(10,16).Value = 0.01363493739935873
$ws.Cells.Item(10,17).Value = 0.3442598793040001
$ws.Cells.Item(10,18).Value = 3.098338913736001
$ws.Cells.Item(10,19).Value = 0.00211940825659222
$ws.Cells.Item(10,20).Value = 0.002119408256592221
